# Update the PSSM score matrix (B2:K21) with supplemental-figure values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,10
$row2[0,0] = -16.52052660637241
$row2[0,1] = 0.05408170260096125
$row2[0,2] = -16.52052660637241
$row2[0,3] = -16.52052660637241
$row2[0,4] = -16.52052660637241
$row2[0,5] = -16.52052660637241
$row2[0,6] = -16.52052660637241
$row2[0,7] = -16.52052660637241
$row2[0,8] = -16.52052660637241
$row2[0,9] = -16.52052660637241
$ws.Range("B2:K2").Value = $row2

$row3 = New-Object "object[,]" 1,10
$row3[0,0] = -16.52052660637241
$row3[0,1] = -16.52052660637241
$row3[0,2] = -16.52052660637241
$row3[0,3] = -16.52052660637241
$row3[0,4] = -16.52052660637241
$row3[0,5] = -16.52052660637241
$row3[0,6] = -16.52052660637241
$row3[0,7] = 1.254549920462841
$row3[0,8] = -16.52052660637241
$row3[0,9] = -16.52052660637241
$ws.Range("B3:K3").Value = $row3

$row4 = New-Object "object[,]" 1,10
$row4[0,0] = -16.52052660637241
$row4[0,1] = -0.05169576992706869
$row4[0,2] = 0.3386320247279529
$row4[0,3] = -16.52052660637241
$row4[0,4] = 3.827904398765399
$row4[0,5] = -16.52052660637241
$row4[0,6] = 1.367104564609802
$row4[0,7] = -16.52052660637241
$row4[0,8] = -16.52052660637241
$row4[0,9] = -16.52052660637241
$ws.Range("B4:K4").Value = $row4

$row5 = New-Object "object[,]" 1,10
$row5[0,0] = -16.52052660637241
$row5[0,1] = 0.3514387744035046
$row5[0,2] = -16.52052660637241
$row5[0,3] = -16.52052660637241
$row5[0,4] = -16.52052660637241
$row5[0,5] = 3.047253116001143
$row5[0,6] = -16.52052660637241
$row5[0,7] = -16.52052660637241
$row5[0,8] = -16.52052660637241
$row5[0,9] = -16.52052660637241
$ws.Range("B5:K5").Value = $row5

$row6 = New-Object "object[,]" 1,10
$row6[0,0] = -16.52052660637241
$row6[0,1] = -16.52052660637241
$row6[0,2] = -16.52052660637241
$row6[0,3] = -16.52052660637241
$row6[0,4] = -16.52052660637241
$row6[0,5] = -16.52052660637241
$row6[0,6] = -16.52052660637241
$row6[0,7] = -16.52052660637241
$row6[0,8] = -16.52052660637241
$row6[0,9] = -16.52052660637241
$ws.Range("B6:K6").Value = $row6

$row7 = New-Object "object[,]" 1,10
$row7[0,0] = 3.199512192744305
$row7[0,1] = -16.52052660637241
$row7[0,2] = -16.52052660637241
$row7[0,3] = -16.52052660637241
$row7[0,4] = -16.52052660637241
$row7[0,5] = -16.52052660637241
$row7[0,6] = -16.52052660637241
$row7[0,7] = -16.52052660637241
$row7[0,8] = -16.52052660637241
$row7[0,9] = -16.52052660637241
$ws.Range("B7:K7").Value = $row7

$row8 = New-Object "object[,]" 1,10
$row8[0,0] = -16.52052660637241
$row8[0,1] = -16.52052660637241
$row8[0,2] = -16.52052660637241
$row8[0,3] = 1.646504351270308
$row8[0,4] = -16.52052660637241
$row8[0,5] = -16.52052660637241
$row8[0,6] = -16.52052660637241
$row8[0,7] = -16.52052660637241
$row8[0,8] = -16.52052660637241
$row8[0,9] = -16.52052660637241
$ws.Range("B8:K8").Value = $row8

$row9 = New-Object "object[,]" 1,10
$row9[0,0] = 3.434738751432035
$row9[0,1] = -16.52052660637241
$row9[0,2] = -16.52052660637241
$row9[0,3] = -16.52052660637241
$row9[0,4] = -16.52052660637241
$row9[0,5] = -16.52052660637241
$row9[0,6] = -16.52052660637241
$row9[0,7] = -16.52052660637241
$row9[0,8] = -16.52052660637241
$row9[0,9] = -16.52052660637241
$ws.Range("B9:K9").Value = $row9

$row10 = New-Object "object[,]" 1,10
$row10[0,0] = -16.52052660637241
$row10[0,1] = -16.52052660637241
$row10[0,2] = -16.52052660637241
$row10[0,3] = -16.52052660637241
$row10[0,4] = -16.52052660637241
$row10[0,5] = -16.52052660637241
$row10[0,6] = -16.52052660637241
$row10[0,7] = 1.594501798819205
$row10[0,8] = -16.52052660637241
$row10[0,9] = 1.583585802634637
$ws.Range("B10:K10").Value = $row10

$row11 = New-Object "object[,]" 1,10
$row11[0,0] = -16.52052660637241
$row11[0,1] = -16.52052660637241
$row11[0,2] = -16.52052660637241
$row11[0,3] = 2.340502731257982
$row11[0,4] = -16.52052660637241
$row11[0,5] = 2.265056766117707
$row11[0,6] = -16.52052660637241
$row11[0,7] = -16.52052660637241
$row11[0,8] = -16.52052660637241
$row11[0,9] = 1.869764125833906
$ws.Range("B11:K11").Value = $row11

$row12 = New-Object "object[,]" 1,10
$row12[0,0] = -16.52052660637241
$row12[0,1] = -16.52052660637241
$row12[0,2] = -16.52052660637241
$row12[0,3] = -16.52052660637241
$row12[0,4] = -16.52052660637241
$row12[0,5] = -16.52052660637241
$row12[0,6] = -16.52052660637241
$row12[0,7] = -16.52052660637241
$row12[0,8] = -16.52052660637241
$row12[0,9] = -16.52052660637241
$ws.Range("B12:K12").Value = $row12

$row13 = New-Object "object[,]" 1,10
$row13[0,0] = -16.52052660637241
$row13[0,1] = -16.52052660637241
$row13[0,2] = -16.52052660637241
$row13[0,3] = 1.891683481059117
$row13[0,4] = -16.52052660637241
$row13[0,5] = -16.52052660637241
$row13[0,6] = -16.52052660637241
$row13[0,7] = -16.52052660637241
$row13[0,8] = 4.321913515936039
$row13[0,9] = 2.139375241433322
$ws.Range("B13:K13").Value = $row13

$row14 = New-Object "object[,]" 1,10
$row14[0,0] = -16.52052660637241
$row14[0,1] = -16.52052660637241
$row14[0,2] = 1.137077953075452
$row14[0,3] = -16.52052660637241
$row14[0,4] = -16.52052660637241
$row14[0,5] = -16.52052660637241
$row14[0,6] = -16.52052660637241
$row14[0,7] = -16.52052660637241
$row14[0,8] = -16.52052660637241
$row14[0,9] = 2.145145876334274
$ws.Range("B14:K14").Value = $row14

$row15 = New-Object "object[,]" 1,10
$row15[0,0] = -16.52052660637241
$row15[0,1] = -16.52052660637241
$row15[0,2] = 0.2690197741971133
$row15[0,3] = -16.52052660637241
$row15[0,4] = -16.52052660637241
$row15[0,5] = -16.52052660637241
$row15[0,6] = -16.52052660637241
$row15[0,7] = -16.52052660637241
$row15[0,8] = -16.52052660637241
$row15[0,9] = -16.52052660637241
$ws.Range("B15:K15").Value = $row15

$row16 = New-Object "object[,]" 1,10
$row16[0,0] = -16.52052660637241
$row16[0,1] = -16.52052660637241
$row16[0,2] = -16.52052660637241
$row16[0,3] = -16.52052660637241
$row16[0,4] = -16.52052660637241
$row16[0,5] = -16.52052660637241
$row16[0,6] = -16.52052660637241
$row16[0,7] = -16.52052660637241
$row16[0,8] = -16.52052660637241
$row16[0,9] = -16.52052660637241
$ws.Range("B16:K16").Value = $row16

$row17 = New-Object "object[,]" 1,10
$row17[0,0] = -16.52052660637241
$row17[0,1] = 0.5085556777670397
$row17[0,2] = 0.1980196477416181
$row17[0,3] = -16.52052660637241
$row17[0,4] = -16.52052660637241
$row17[0,5] = -16.52052660637241
$row17[0,6] = 1.708954106687743
$row17[0,7] = 0.1993568261384911
$row17[0,8] = -16.52052660637241
$row17[0,9] = -16.52052660637241
$ws.Range("B17:K17").Value = $row17

$row18 = New-Object "object[,]" 1,10
$row18[0,0] = -16.52052660637241
$row18[0,1] = -16.52052660637241
$row18[0,2] = -16.52052660637241
$row18[0,3] = -16.52052660637241
$row18[0,4] = -16.52052660637241
$row18[0,5] = -16.52052660637241
$row18[0,6] = 1.728628521784999
$row18[0,7] = 0.8981248775407245
$row18[0,8] = -16.52052660637241
$row18[0,9] = -16.52052660637241
$ws.Range("B18:K18").Value = $row18

$row19 = New-Object "object[,]" 1,10
$row19[0,0] = -16.52052660637241
$row19[0,1] = -16.52052660637241
$row19[0,2] = 2.851634903787438
$row19[0,3] = -16.52052660637241
$row19[0,4] = -16.52052660637241
$row19[0,5] = -16.52052660637241
$row19[0,6] = 2.098177920373745
$row19[0,7] = 2.225374171258364
$row19[0,8] = -16.52052660637241
$row19[0,9] = -16.52052660637241
$ws.Range("B19:K19").Value = $row19

$row20 = New-Object "object[,]" 1,10
$row20[0,0] = -16.52052660637241
$row20[0,1] = 3.167349660178635
$row20[0,2] = 2.800258520733087
$row20[0,3] = -16.52052660637241
$row20[0,4] = 2.535796225660373
$row20[0,5] = -16.52052660637241
$row20[0,6] = 1.609163922412888
$row20[0,7] = 2.787823534780015
$row20[0,8] = -16.52052660637241
$row20[0,9] = 2.175968952113783
$ws.Range("B20:K20").Value = $row20

$row21 = New-Object "object[,]" 1,10
$row21[0,0] = -16.52052660637241
$row21[0,1] = 2.658655997859987
$row21[0,2] = -16.52052660637241
$row21[0,3] = 3.01678470532227
$row21[0,4] = -16.52052660637241
$row21[0,5] = 2.792172473679509
$row21[0,6] = 1.809221564142468
$row21[0,7] = -16.52052660637241
$row21[0,8] = -16.52052660637241
$row21[0,9] = -16.52052660637241
$ws.Range("B21:K21").Value = $row21
